$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel
# (e.g. "1.00", "6.89") are forced to remain plain text by temporarily applying
# a text number format, then restoring the default "Normal" style so no visible
# formatting change is left behind.
$textCells = @("D4", "D5", "D6", "D8", "D12", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D30", "D32", "D33", "D34", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '60.890.59'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.643.37'
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '577.29'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '143.94'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '0.377'
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").Value = '3.113.64'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '26.39'
$ws.Range("E14").Value = '  +6.33%  '
$ws.Range("D15").Value = '60.834.93'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("E16").Value = '  +1.50%  '
$ws.Range("D17").Value = '2.661.91'
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").Value = '11.63'
$ws.Range("E18").Value = '  +2.05%  '
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").Value = '352.09'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '6.89'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '0.526'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '63.97'
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '8.40'
$ws.Range("E27").Value = '  +5.98%  '
$ws.Range("D28").Value = '2.02'
$ws.Range("E28").Value = '  +8.63%  '
$ws.Range("D29").Value = '0.0₃0808'
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("D30").Value = '6.77'
$ws.Range("E30").Value = '  +6.11%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '166.65'
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("D33").Value = '19.97'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").Value = '4.56'
$ws.Range("E34").Value = '  +7.21%  '
$ws.Range("E35").Value = '  +6.33%  '
$ws.Range("E36").Value = '  +7.16%  '
$ws.Range("E37").Value = '  +3.90%  '
$ws.Range("D38").Value = '341.78'
$ws.Range("E38").Value = '  +8.24%  '
$ws.Range("D39").Value = '4.13'
$ws.Range("E39").Value = '  +5.95%  '
$ws.Range("D40").Value = '0.910'
$ws.Range("E40").Value = '  +7.56%  '
$ws.Range("D41").Value = '38.32'
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").Value = '138.29'
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("E43").Value = '  +4.64%  '
$ws.Range("D44").Value = '0.0574'
$ws.Range("E44").Value = '  +3.77%  '
$ws.Range("D45").Value = '0.624'
$ws.Range("E45").Value = '  +2.68%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0252'
$ws.Range("E46").Value = '  +3.45%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '20.83'
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '20.31'
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("D49").Value = '0.0996'
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").Value = '2.089.33'
$ws.Range("E51").Value = '  +2.47%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
